$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.ClearContents()

# Header row
$ws.Range("A1").Value = "Sending cluster"
$ws.Range("B1").Value = "Ligand symbol"
$ws.Range("C1").Value = "Receptor symbol"
$ws.Range("D1").Value = "Target cluster"
$ws.Range("E1").Value = "Ligand-expressing cells"
$ws.Range("F1").Value = "Ligand detection rate"
$ws.Range("G1").Value = "Ligand average expression value"
$ws.Range("H1").Value = "Ligand total expression value"
$ws.Range("I1").Value = "Ligand derived specificity of average expression value"
$ws.Range("J1").Value = "Ligand derived specificity of total expression value"
$ws.Range("K1").Value = "Receptor-expressing cells"
$ws.Range("L1").Value = "Receptor detection rate"
$ws.Range("M1").Value = "Receptor average expression value"
$ws.Range("N1").Value = "Receptor total expression value"
$ws.Range("O1").Value = "Receptor derived specificity of average expression value"
$ws.Range("P1").Value = "Receptor derived specificity of total expression value"
$ws.Range("Q1").Value = "Edge average expression weight"
$ws.Range("R1").Value = "Edge total expression weight"
$ws.Range("S1").Value = "Edge average expression derived specificity"
$ws.Range("T1").Value = "Edge total expression derived specificity"

# Data rows, written column-major to control shared-string ordering
$ws.Range("A2").Value = "FAPs"
$ws.Range("A3").Value = "FAPs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("A6").Value = "sCs"
$ws.Range("A7").Value = "sCs"
$ws.Range("A8").Value = "sCs"
$ws.Range("A9").Value = "sCs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C2").Value = "Lrp5"
$ws.Range("C3").Value = "Lrp5"
$ws.Range("C4").Value = "Lrp5"
$ws.Range("C5").Value = "Lrp5"
$ws.Range("C6").Value = "Lrp5"
$ws.Range("C7").Value = "Lrp5"
$ws.Range("C8").Value = "Lrp5"
$ws.Range("C9").Value = "Lrp5"
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "M2"
$ws.Range("D5").Value = "sCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("D8").Value = "M2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 2
$ws.Range("E9").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G2").Value = 10.43424333333333
$ws.Range("G3").Value = 10.43424333333333
$ws.Range("G4").Value = 10.43424333333333
$ws.Range("G5").Value = 10.43424333333333
$ws.Range("G6").Value = 0.310422
$ws.Range("G7").Value = 0.310422
$ws.Range("G8").Value = 0.310422
$ws.Range("G9").Value = 0.310422
$ws.Range("H2").Value = 31.30273
$ws.Range("H3").Value = 31.30273
$ws.Range("H4").Value = 31.30273
$ws.Range("H5").Value = 31.30273
$ws.Range("H6").Value = 0.9312659999999999
$ws.Range("H7").Value = 0.9312659999999999
$ws.Range("H8").Value = 0.9312659999999999
$ws.Range("H9").Value = 0.9312659999999999
$ws.Range("I2").Value = 0.9711091978791583
$ws.Range("I3").Value = 0.9711091978791583
$ws.Range("I4").Value = 0.9711091978791583
$ws.Range("I5").Value = 0.9711091978791583
$ws.Range("I6").Value = 0.02889080212084161
$ws.Range("I7").Value = 0.02889080212084161
$ws.Range("I8").Value = 0.02889080212084161
$ws.Range("I9").Value = 0.02889080212084161
$ws.Range("J2").Value = 0.9711091978791584
$ws.Range("J3").Value = 0.9711091978791584
$ws.Range("J4").Value = 0.9711091978791584
$ws.Range("J5").Value = 0.9711091978791584
$ws.Range("J6").Value = 0.02889080212084161
$ws.Range("J7").Value = 0.02889080212084161
$ws.Range("J8").Value = 0.02889080212084161
$ws.Range("J9").Value = 0.02889080212084161
$ws.Range("K2").Value = 3
$ws.Range("K3").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("K8").Value = 3
$ws.Range("K9").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("L9").Value = 1
$ws.Range("M2").Value = 13.523597
$ws.Range("M3").Value = 11.59690133333333
$ws.Range("M4").Value = 4.574811666666666
$ws.Range("M5").Value = 9.452519666666667
$ws.Range("M6").Value = 13.523597
$ws.Range("M7").Value = 11.59690133333333
$ws.Range("M8").Value = 4.574811666666666
$ws.Range("M9").Value = 9.452519666666667
$ws.Range("N2").Value = 40.570791
$ws.Range("N3").Value = 34.790704
$ws.Range("N4").Value = 13.724435
$ws.Range("N5").Value = 28.357559
$ws.Range("N6").Value = 40.570791
$ws.Range("N7").Value = 34.790704
$ws.Range("N8").Value = 13.724435
$ws.Range("N9").Value = 28.357559
$ws.Range("O2").Value = 0.3454494697445509
$ws.Range("O3").Value = 0.2962335698320407
$ws.Range("O4").Value = 0.1168599052775075
$ws.Range("O5").Value = 0.2414570551459009
$ws.Range("O6").Value = 0.3454494697445509
$ws.Range("O7").Value = 0.2962335698320407
$ws.Range("O8").Value = 0.1168599052775075
$ws.Range("O9").Value = 0.2414570551459009
$ws.Range("P2").Value = 0.3454494697445509
$ws.Range("P3").Value = 0.2962335698320406
$ws.Range("P4").Value = 0.1168599052775075
$ws.Range("P5").Value = 0.2414570551459009
$ws.Range("P6").Value = 0.3454494697445509
$ws.Range("P7").Value = 0.2962335698320406
$ws.Range("P8").Value = 0.1168599052775075
$ws.Range("P9").Value = 0.2414570551459009
$ws.Range("Q2").Value = 141.1085018399367
$ws.Range("Q3").Value = 121.0048904246578
$ws.Range("Q4").Value = 47.73469813417221
$ws.Range("Q5").Value = 98.62989031511889
$ws.Range("Q6").Value = 4.198022027934
$ws.Range("Q7").Value = 3.599933305696
$ws.Range("Q8").Value = 1.42012218719
$ws.Range("Q9").Value = 2.934270059966
$ws.Range("R2").Value = 1269.97651655943
$ws.Range("R3").Value = 1089.04401382192
$ws.Range("R4").Value = 429.6122832075499
$ws.Range("R5").Value = 887.66901283607
$ws.Range("R6").Value = 37.78219825140599
$ws.Range("R7").Value = 32.39939975126399
$ws.Range("R8").Value = 12.78109968471
$ws.Range("R9").Value = 26.408430539694
$ws.Range("S2").Value = 0.3354691574714114
$ws.Range("S3").Value = 0.2876751443844727
$ws.Range("S4").Value = 0.1134837288782748
$ws.Range("S5").Value = 0.2344811671449995
$ws.Range("S6").Value = 0.009980312273139481
$ws.Range("S7").Value = 0.008558425447568001
$ws.Range("S8").Value = 0.003376176399232764
$ws.Range("S9").Value = 0.006975888000901362
$ws.Range("T2").Value = 0.3354691574714114
$ws.Range("T3").Value = 0.2876751443844726
$ws.Range("T4").Value = 0.1134837288782748
$ws.Range("T5").Value = 0.2344811671449995
$ws.Range("T6").Value = 0.009980312273139481
$ws.Range("T7").Value = 0.008558425447568001
$ws.Range("T8").Value = 0.003376176399232765
$ws.Range("T9").Value = 0.006975888000901363
